# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28, ~15:17-15:18) to the
# PIR, Humidity and Temperature sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

# NOTE: this runtime's PowerShell-style interpreter only binds positional
# function arguments, so the helper below is called as
# `Append-LogRows $SheetName $StartRow $Rows $TextCols` (no `-Name value`
# switches).
#
# $TextCols lists the 1-based column numbers whose values look like a
# number/date/percentage to Excel (e.g. "2026-01-28" or "88.3%") and must be
# forced to plain text so they land in the sheet as literal strings -- the
# same way the rest of this log is stored -- instead of being auto-converted
# to a date serial / numeric value.
function Append-LogRows($SheetName, $StartRow, $Rows, $TextCols) {
    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    foreach ($col in $TextCols) {
        $colLetter = [char](64 + $col)
        $ws.Range("$colLetter$StartRow" + ":$colLetter$endRow").NumberFormat = "@"
    }

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $values = $Rows[$i]
        $ws.Cells.Item($r, 1).Value = $values[0]
        $ws.Cells.Item($r, 2).Value = $values[1]
        $ws.Cells.Item($r, 3).Value = $values[2]
        $ws.Cells.Item($r, 4).Value = $values[3]
        $ws.Cells.Item($r, 5).Value = $values[4]
        $ws.Cells.Item($r, 6).Value = $values[5]
    }
}

# ---------------------------------------------------------------------
# PIR sheet: rows 307-318 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------------
$pirData = @(
    @("2026-01-28","15:17:41","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:17:44","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:17:50","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:17:54","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:17:59","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:04","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:10","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:14","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:19","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:24","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:29","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:18:34","15:00","Bathroom","No Motion","Inactive")
)
Append-LogRows "PIR" 307 $pirData @(1)

# ---------------------------------------------------------------------
# Humidity sheet: rows 291-303
# ---------------------------------------------------------------------
$humidityData = @(
    @("2026-01-28","15:17:39","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:17:42","15:00","Bathroom","87.4%","Active"),
    @("2026-01-28","15:17:45","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:17:48","15:00","Bathroom","88.4%","Active"),
    @("2026-01-28","15:17:52","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:17:56","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:18:00","15:00","Bathroom","87.4%","Active"),
    @("2026-01-28","15:18:05","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:18:08","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:18:12","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:18:21","15:00","Bathroom","87.3%","Active"),
    @("2026-01-28","15:18:33","15:00","Bathroom","87.3%","Active"),
    @("2026-01-28","15:18:37","15:00","Bathroom","88.3%","Active")
)
Append-LogRows "Humidity" 291 $humidityData @(1,5)

# ---------------------------------------------------------------------
# Temperature sheet: rows 291-303
# ---------------------------------------------------------------------
$temperatureData = @(
    @("2026-01-28","15:17:40","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:17:43","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:17:45","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:17:49","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:17:53","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:17:57","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:01","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:06","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:09","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:13","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:21","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:33","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:18:37","15:00","Bathroom","22.9C","Active")
)
Append-LogRows "Temperature" 291 $temperatureData @(1)

Write-Output "Appended $($pirData.Count) PIR rows, $($humidityData.Count) Humidity rows, $($temperatureData.Count) Temperature rows."
